# Swap the presentation's theme colour scheme from the custom
# "Integral" palette to the stock PowerPoint "Office Theme" palette
# (ppt/theme/theme1.xml), matching the tonal reset captured in the
# commit's OOXML diff (the deck's slide master picks up the plain
# default Office colours in place of the bespoke Integral green/teal
# set).
#
# PowerPoint's ThemeColorScheme exposes the 12 theme colour slots in
# the same fixed order as <a:clrScheme>: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. COM RGB values are encoded as decimal
# BGR (0xBBGGRR), so each hex target below is byte-swapped before
# assignment.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# Ordered target sRGB values (as they appear in <a:clrScheme>) for the
# default Office theme.
$officeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $bgrValue = $b * 65536 + $g * 256 + $r
    $cs.Colors($i).RGB = $bgrValue
}
